$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the current row 454, pushing every
# subsequent row (454..564) down by one (to 455..565). Excel copies the
# formatting (incl. the date number format on column D) from the row
# above automatically on insert, matching the existing sheet's style.
$ws.Rows.Item(454).Insert()

# Populate the newly inserted row 454 with the new weekly record.
$ws.Cells.Item(454, 1).Value = 7
$ws.Cells.Item(454, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(454, 3).Value = 'Ñuble'
$ws.Cells.Item(454, 4).Value = (Get-Date -Year 2023 -Month 7 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(454, 5).Value = 16
$ws.Cells.Item(454, 6).Value = 100112023
$ws.Cells.Item(454, 7).Value = 'Brócoli'
$ws.Cells.Item(454, 8).Value = 'Sin especificar'
$ws.Cells.Item(454, 9).Value = 'Primera'
$ws.Cells.Item(454, 10).Value = 150
$ws.Cells.Item(454, 11).Value = 1000
$ws.Cells.Item(454, 12).Value = 1000
$ws.Cells.Item(454, 13).Value = 1000
$ws.Cells.Item(454, 14).Value = '$/unidad'
$ws.Cells.Item(454, 15).Value = 'Región del Maule'
$ws.Cells.Item(454, 16).Value = 1000
$ws.Cells.Item(454, 17).Value = 1
$ws.Cells.Item(454, 18).Value = 'Hortaliza'
